$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 98

# A98: date value, reuse the date/time number format used by the cells above (style index 1)
$ws.Cells.Item($row, 1).Value = 45449.2916666667
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122) # xlPasteFormats

# B98-F98: plain numeric values
$ws.Cells.Item($row, 2).Value = 0
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 1
$ws.Cells.Item($row, 6).Value = 1

# G98: text "1" (would otherwise be auto-converted to a number) -> use quote-prefix
# to force text, then restore default style so the cell keeps default formatting.
$ws.Cells.Item($row, 7).Value = "'1"
$ws.Cells.Item($row, 7).Style = "Normal"

# H98: ticker text
$ws.Cells.Item($row, 8).Value = "YKY.MI"
